$d = $word.ActiveDocument

# Replace the two title words "Identity" and "Provider" with
# "Decentralised" and "Identifiers" respectively, matching the
# whole word only so we don't touch other occurrences.

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Identity", $true, $true, $false, $false, $false, $true, 1, $false, "Decentralised", 2)

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("Provider", $true, $true, $false, $false, $false, $true, 1, $false, "Identifiers", 2)
